$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 583.3333
$ws.Range("I18").Value = 583.3333
$ws.Range("K18").Value = 583.3333
$ws.Range("M18").Value = -299.3333

$ws.Range("H80").Value = 724.1429000000001
$ws.Range("I80").Value = 456.75
$ws.Range("K80").Value = 1370.25
$ws.Range("M80").Value = -372.25

$ws.Range("H83").Value = 724.1429000000001
$ws.Range("I83").Value = 456.75
$ws.Range("K83").Value = 4110.75
$ws.Range("M83").Value = 881.25

$ws.Range("H135").Value = 2599.2
$ws.Range("I135").Value = 2027.4286
$ws.Range("K135").Value = 18246.8574
$ws.Range("M135").Value = -15711.8574

$ws.Range("H138").Value = 15670.685
$ws.Range("J138").Value = 41689.652
$ws.Range("L138").Value = 125068.956
$ws.Range("N138").Value = -135348.956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5640.9585
$ws.Range("I61").Value = 934.3889
$ws.Range("K61").Value = 934.3889
$ws.Range("M61").Value = -722.3889

$ws.Range("H74").Value = 291654.56
$ws.Range("I74").Value = 400717.66
$ws.Range("K74").Value = 400717.66
$ws.Range("M74").Value = -399843.66

$ws.Range("H77").Value = 291654.56
$ws.Range("I77").Value = 400717.66
$ws.Range("K77").Value = 2003588.3
$ws.Range("M77").Value = -1999220.3

$ws.Range("H132").Value = 1396.3334
$ws.Range("I132").Value = 1113.0938
$ws.Range("J132").Value = 3662.25
$ws.Range("K132").Value = 3339.2814
$ws.Range("L132").Value = 10986.75
$ws.Range("M132").Value = -809.2814000000003
$ws.Range("N132").Value = -16046.75

$ws.Range("H136").Value = 5640.9585
$ws.Range("I136").Value = 934.3889
$ws.Range("K136").Value = 2803.1667
$ws.Range("M136").Value = -253.1667000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 21602.25
$ws.Range("I82").Value = 14903.375
$ws.Range("K82").Value = 14903.375
$ws.Range("M82").Value = -14520.375

$ws.Range("H85").Value = 21602.25
$ws.Range("I85").Value = 14903.375
$ws.Range("K85").Value = 14903.375
$ws.Range("M85").Value = -13577.375

$ws.Range("H94").Value = 1225.6666
$ws.Range("I94").Value = 1039.6666
$ws.Range("K94").Value = 1039.6666
$ws.Range("M94").Value = -588.6666

$ws.Range("H105").Value = 2736.2856
$ws.Range("I105").Value = 1641.6428
$ws.Range("K105").Value = 1641.6428
$ws.Range("M105").Value = 105.3571999999999

$ws.Range("H107").Value = 5217.231
$ws.Range("I107").Value = 4168.45
$ws.Range("J107").Value = 8713.166999999999
$ws.Range("K107").Value = 4168.45
$ws.Range("L107").Value = 8713.166999999999
$ws.Range("M107").Value = -2248.45
$ws.Range("N107").Value = -12553.167

$ws.Range("H134").Value = 8928.172
$ws.Range("I134").Value = 9672.138000000001
$ws.Range("K134").Value = 29016.414
$ws.Range("M134").Value = -26481.414

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13789.852
$ws.Range("I58").Value = 1382.1904
$ws.Range("K58").Value = 1382.1904
$ws.Range("M58").Value = -1179.1904

$ws.Range("H62").Value = 7455.4165
$ws.Range("J62").Value = 3742.5
$ws.Range("L62").Value = 3742.5
$ws.Range("N62").Value = -4990.5

$ws.Range("H65").Value = 7455.4165
$ws.Range("J65").Value = 3742.5
$ws.Range("L65").Value = 18712.5
$ws.Range("N65").Value = -24952.5

$ws.Range("H107").Value = 680.375
$ws.Range("I107").Value = 625.6923
$ws.Range("J107").Value = 745
$ws.Range("K107").Value = 625.6923
$ws.Range("L107").Value = 745
$ws.Range("M107").Value = 1294.3077
$ws.Range("N107").Value = -4585

$ws.Range("H122").Value = 2410.5
$ws.Range("I122").Value = 2301.2
$ws.Range("K122").Value = 6903.599999999999
$ws.Range("M122").Value = -4453.599999999999

$ws.Range("H132").Value = 51459.85
$ws.Range("I132").Value = 56622.055
$ws.Range("K132").Value = 169866.165
$ws.Range("M132").Value = -167336.165

$ws.Range("H136").Value = 13789.852
$ws.Range("I136").Value = 1382.1904
$ws.Range("K136").Value = 4146.5712
$ws.Range("M136").Value = -1596.5712

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 82127620
$ws.Range("I4").Value = 108552060
$ws.Range("K4").Value = 325656180
$ws.Range("M4").Value = -325656068

$ws.Range("H121").Value = 60512.652
$ws.Range("I121").Value = 133694.62
$ws.Range("J121").Value = 21482.268
$ws.Range("K121").Value = 401083.86
$ws.Range("L121").Value = 64446.804
$ws.Range("M121").Value = -399773.86
$ws.Range("N121").Value = -67066.804

$ws.Range("H132").Value = 1025.7273
$ws.Range("I132").Value = 737.4
$ws.Range("J132").Value = 1266
$ws.Range("K132").Value = 6636.599999999999
$ws.Range("L132").Value = 11394
$ws.Range("M132").Value = -4106.599999999999
$ws.Range("N132").Value = -16454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 10000
$ws.Range("I12").Value = 10000
$ws.Range("K12").Value = 10000
$ws.Range("M12").Value = -9860

$ws.Range("H102").Value = 2614.7083
$ws.Range("I102").Value = 1640.7142
$ws.Range("J102").Value = 3978.3
$ws.Range("K102").Value = 1640.7142
$ws.Range("L102").Value = 3978.3
$ws.Range("M102").Value = -18.71419999999989
$ws.Range("N102").Value = -7222.3

$ws.Range("H126").Value = 3142.8333
$ws.Range("I126").Value = 1871.4
$ws.Range("K126").Value = 5614.200000000001
$ws.Range("M126").Value = -3144.200000000001

$ws.Range("H132").Value = 3876.2
$ws.Range("I132").Value = 3857.75
$ws.Range("J132").Value = 3950
$ws.Range("K132").Value = 11573.25
$ws.Range("L132").Value = 11850
$ws.Range("M132").Value = -9043.25
$ws.Range("N132").Value = -16910

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2635.125
$ws.Range("I132").Value = 2409.16
$ws.Range("K132").Value = 7227.48
$ws.Range("M132").Value = -4697.48

$ws.Range("H136").Value = 5375
$ws.Range("J136").Value = 5809.143
$ws.Range("L136").Value = 17427.429
$ws.Range("N136").Value = -22527.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4577.6665
$ws.Range("I62").Value = 3733
$ws.Range("K62").Value = 3733
$ws.Range("M62").Value = -3109

$ws.Range("H65").Value = 4577.6665
$ws.Range("I65").Value = 3733
$ws.Range("K65").Value = 18665
$ws.Range("M65").Value = -15545

$ws.Range("H122").Value = 56566.117
$ws.Range("J122").Value = 3068.8572
$ws.Range("L122").Value = 9206.571599999999
$ws.Range("N122").Value = -14106.5716

$ws.Range("H126").Value = 4950
$ws.Range("I126").Value = 4900
$ws.Range("K126").Value = 14700
$ws.Range("M126").Value = -12230

$ws.Range("H132").Value = 883.4737
$ws.Range("J132").Value = 2483.3333
$ws.Range("L132").Value = 7449.999899999999
$ws.Range("N132").Value = -12509.9999

$ws.Range("H136").Value = 39893.355
$ws.Range("I136").Value = 59335.445
$ws.Range("K136").Value = 178006.335
$ws.Range("M136").Value = -175456.335
